$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("000 - Printed Parts")

# Insert two new rows at the top of the printed-parts list (row 87) for the
# new PN120 / PN121 entries, pushing the existing rows down by two.
$ws.Rows("87:88").Insert()

# Row 87: PN120 - Exterior PTFE Coupler
$ws.Range("A87").Value = "PN120"
$ws.Range("B87").Value = "Exterior"
$ws.Range("C87").Value = "Frame"
$ws.Range("D87").Value = "R"
$ws.Range("E87").Value = "Exterior PTFE Coupler"
$ws.Range("F87").Value = "ABS"
$ws.Range("G87").Value = 1
$ws.Range("I87").Value = "120 - Exterior - Frame - Exterior PTFE Coupler.stl"

# Row 88: PN121 - Exhaust Hose Coupler
$ws.Range("A88").Value = "PN121"
$ws.Range("B88").Value = "Exterior"
$ws.Range("C88").Value = "Frame"
$ws.Range("D88").Value = "L"
$ws.Range("E88").Value = "Exhaust Hose Coupler"
$ws.Range("F88").Value = "ABS"
$ws.Range("G88").Value = 1
$ws.Range("I88").Value = "121 - Exterior - Frame - Exhaust Hose Coupler.stl"

# Grow the table (and its autofilter) to include the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I99"))

# Reposition the view roughly where the author last left it.
$ws.Application.ActiveWindow.ScrollRow = 82
$ws.Range("G92").Select()
